$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2020" column (Q), copying the formatting
# from the preceding "2019" column (P) so the new cells pick up the same
# fonts/borders/number formats already used by the rest of the table.
$ws.Range("P2:P5").Copy()
$ws.Range("Q2:Q5").PasteSpecial(-4122)

# Fill in the new year's data.
$ws.Range("Q3").Value = 2020
$ws.Range("Q4").Value = 197792
$ws.Range("Q5").Value = 2.9794303052841493

# Restore the cursor/selection that was active when the workbook was saved.
$ws.Range("G27").Select()
